# Natmi following Dr Hou advice
# Update LR-pairs table: add "ECs" as a third sending cluster alongside
# "FAPs" and "sCs", expanding the 2x3 sender/target grid to a 3x3 grid
# (rows 2-7 updated in place, rows 8-10 added).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tnc"
$ws.Cells.Item(2, 3).Value = "Itgb1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 1.442371333333333
$ws.Cells.Item(2, 8).Value = 4.327114
$ws.Cells.Item(2, 9).Value = 0.03522044016446201
$ws.Cells.Item(2, 10).Value = 0.03522044016446201
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 112.513392
$ws.Cells.Item(2, 14).Value = 337.540176
$ws.Cells.Item(2, 15).Value = 0.3275312977368564
$ws.Cells.Item(2, 16).Value = 0.3275312977368564
$ws.Cells.Item(2, 17).Value = 162.286091236896
$ws.Cells.Item(2, 18).Value = 1460.574821132064
$ws.Cells.Item(2, 19).Value = 0.01153579647392954
$ws.Cells.Item(2, 20).Value = 0.01153579647392954
# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tnc"
$ws.Cells.Item(3, 3).Value = "Itgb1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 1.442371333333333
$ws.Cells.Item(3, 8).Value = 4.327114
$ws.Cells.Item(3, 9).Value = 0.03522044016446201
$ws.Cells.Item(3, 10).Value = 0.03522044016446201
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 106.314466
$ws.Cells.Item(3, 14).Value = 318.943398
$ws.Cells.Item(3, 15).Value = 0.3094859589441663
$ws.Cells.Item(3, 16).Value = 0.3094859589441664
$ws.Cells.Item(3, 17).Value = 153.3449380770413
$ws.Cells.Item(3, 18).Value = 1380.104442693372
$ws.Cells.Item(3, 19).Value = 0.01090023169873416
$ws.Cells.Item(3, 20).Value = 0.01090023169873416
# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tnc"
$ws.Cells.Item(4, 3).Value = "Itgb1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 1.442371333333333
$ws.Cells.Item(4, 8).Value = 4.327114
$ws.Cells.Item(4, 9).Value = 0.03522044016446201
$ws.Cells.Item(4, 10).Value = 0.03522044016446201
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 124.6916553333333
$ws.Cells.Item(4, 14).Value = 374.074966
$ws.Cells.Item(4, 15).Value = 0.3629827433189773
$ws.Cells.Item(4, 16).Value = 0.3629827433189773
$ws.Cells.Item(4, 17).Value = 179.8516691586804
$ws.Cells.Item(4, 18).Value = 1618.665022428124
$ws.Cells.Item(4, 19).Value = 0.01278441199179831
$ws.Cells.Item(4, 20).Value = 0.01278441199179831
# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Tnc"
$ws.Cells.Item(5, 3).Value = "Itgb1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 11.331397
$ws.Cells.Item(5, 8).Value = 33.994191
$ws.Cells.Item(5, 9).Value = 0.2766948987373093
$ws.Cells.Item(5, 10).Value = 0.2766948987373092
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 112.513392
$ws.Cells.Item(5, 14).Value = 337.540176
$ws.Cells.Item(5, 15).Value = 0.3275312977368564
$ws.Cells.Item(5, 16).Value = 0.3275312977368564
$ws.Cells.Item(5, 17).Value = 1274.933912568624
$ws.Cells.Item(5, 18).Value = 11474.40521311762
$ws.Cells.Item(5, 19).Value = 0.09062623926059897
$ws.Cells.Item(5, 20).Value = 0.09062623926059896
# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Tnc"
$ws.Cells.Item(6, 3).Value = "Itgb1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 11.331397
$ws.Cells.Item(6, 8).Value = 33.994191
$ws.Cells.Item(6, 9).Value = 0.2766948987373093
$ws.Cells.Item(6, 10).Value = 0.2766948987373092
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 106.314466
$ws.Cells.Item(6, 14).Value = 318.943398
$ws.Cells.Item(6, 15).Value = 0.3094859589441663
$ws.Cells.Item(6, 16).Value = 0.3094859589441664
$ws.Cells.Item(6, 17).Value = 1204.691421089002
$ws.Cells.Item(6, 18).Value = 10842.22278980102
$ws.Cells.Item(6, 19).Value = 0.08563318607067516
$ws.Cells.Item(6, 20).Value = 0.08563318607067516
# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Tnc"
$ws.Cells.Item(7, 3).Value = "Itgb1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 11.331397
$ws.Cells.Item(7, 8).Value = 33.994191
$ws.Cells.Item(7, 9).Value = 0.2766948987373093
$ws.Cells.Item(7, 10).Value = 0.2766948987373092
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 124.6916553333333
$ws.Cells.Item(7, 14).Value = 374.074966
$ws.Cells.Item(7, 15).Value = 0.3629827433189773
$ws.Cells.Item(7, 16).Value = 0.3629827433189773
$ws.Cells.Item(7, 17).Value = 1412.930649169168
$ws.Cells.Item(7, 18).Value = 12716.37584252251
$ws.Cells.Item(7, 19).Value = 0.1004354734060351
$ws.Cells.Item(7, 20).Value = 0.1004354734060351
# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Tnc"
$ws.Cells.Item(8, 3).Value = "Itgb1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 28.17890933333333
$ws.Cells.Item(8, 8).Value = 84.536728
$ws.Cells.Item(8, 9).Value = 0.6880846610982287
$ws.Cells.Item(8, 10).Value = 0.6880846610982286
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 112.513392
$ws.Cells.Item(8, 14).Value = 337.540176
$ws.Cells.Item(8, 15).Value = 0.3275312977368564
$ws.Cells.Item(8, 16).Value = 0.3275312977368564
$ws.Cells.Item(8, 17).Value = 3170.504671953792
$ws.Cells.Item(8, 18).Value = 28534.54204758412
$ws.Cells.Item(8, 19).Value = 0.2253692620023278
$ws.Cells.Item(8, 20).Value = 0.2253692620023278
# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Tnc"
$ws.Cells.Item(9, 3).Value = "Itgb1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 28.17890933333333
$ws.Cells.Item(9, 8).Value = 84.536728
$ws.Cells.Item(9, 9).Value = 0.6880846610982287
$ws.Cells.Item(9, 10).Value = 0.6880846610982286
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 106.314466
$ws.Cells.Item(9, 14).Value = 318.943398
$ws.Cells.Item(9, 15).Value = 0.3094859589441663
$ws.Cells.Item(9, 16).Value = 0.3094859589441664
$ws.Cells.Item(9, 17).Value = 2995.825698235749
$ws.Cells.Item(9, 18).Value = 26962.43128412174
$ws.Cells.Item(9, 19).Value = 0.212952541174757
$ws.Cells.Item(9, 20).Value = 0.212952541174757
# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Tnc"
$ws.Cells.Item(10, 3).Value = "Itgb1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 28.17890933333333
$ws.Cells.Item(10, 8).Value = 84.536728
$ws.Cells.Item(10, 9).Value = 0.6880846610982287
$ws.Cells.Item(10, 10).Value = 0.6880846610982286
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 124.6916553333333
$ws.Cells.Item(10, 14).Value = 374.074966
$ws.Cells.Item(10, 15).Value = 0.3629827433189773
$ws.Cells.Item(10, 16).Value = 0.3629827433189773
$ws.Cells.Item(10, 17).Value = 3513.67485026125
$ws.Cells.Item(10, 18).Value = 31623.07365235125
$ws.Cells.Item(10, 19).Value = 0.2497628579211438
$ws.Cells.Item(10, 20).Value = 0.2497628579211438
